$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-06-02 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-03 Monday", 2)

$d.Content.Find.Execute("78×82=", $true, $false, $false, $false, $false, $true, 1, $false, "86×62=", 2)
$d.Content.Find.Execute("78×51=", $true, $false, $false, $false, $false, $true, 1, $false, "60×55=", 2)
$d.Content.Find.Execute("71×22=", $true, $false, $false, $false, $false, $true, 1, $false, "40×84=", 2)
$d.Content.Find.Execute("29×18=", $true, $false, $false, $false, $false, $true, 1, $false, "19×70=", 2)
$d.Content.Find.Execute("79×27=", $true, $false, $false, $false, $false, $true, 1, $false, "30×58=", 2)
$d.Content.Find.Execute("99×99=", $true, $false, $false, $false, $false, $true, 1, $false, "47×46=", 2)
$d.Content.Find.Execute("63×16=", $true, $false, $false, $false, $false, $true, 1, $false, "20×72=", 2)
$d.Content.Find.Execute("65×28=", $true, $false, $false, $false, $false, $true, 1, $false, "88×46=", 2)
$d.Content.Find.Execute("18×67=", $true, $false, $false, $false, $false, $true, 1, $false, "62×17=", 2)
$d.Content.Find.Execute("81×17=", $true, $false, $false, $false, $false, $true, 1, $false, "21×82=", 2)
$d.Content.Find.Execute("50×48=", $true, $false, $false, $false, $false, $true, 1, $false, "14×19=", 2)
$d.Content.Find.Execute("34×34=", $true, $false, $false, $false, $false, $true, 1, $false, "31×71=", 2)
$d.Content.Find.Execute("74×15=", $true, $false, $false, $false, $false, $true, 1, $false, "25×80=", 2)
$d.Content.Find.Execute("33×56=", $true, $false, $false, $false, $false, $true, 1, $false, "58×29=", 2)
$d.Content.Find.Execute("75×97=", $true, $false, $false, $false, $false, $true, 1, $false, "57×99=", 2)
$d.Content.Find.Execute("78×30=", $true, $false, $false, $false, $false, $true, 1, $false, "84×21=", 2)
$d.Content.Find.Execute("76×99=", $true, $false, $false, $false, $false, $true, 1, $false, "31×60=", 2)
$d.Content.Find.Execute("36×60=", $true, $false, $false, $false, $false, $true, 1, $false, "31×79=", 2)
$d.Content.Find.Execute("32×53=", $true, $false, $false, $false, $false, $true, 1, $false, "20×54=", 2)
$d.Content.Find.Execute("63×43=", $true, $false, $false, $false, $false, $true, 1, $false, "38×99=", 2)
$d.Content.Find.Execute("50×78=", $true, $false, $false, $false, $false, $true, 1, $false, "33×32=", 2)
$d.Content.Find.Execute("90×22=", $true, $false, $false, $false, $false, $true, 1, $false, "26×56=", 2)
$d.Content.Find.Execute("12×90=", $true, $false, $false, $false, $false, $true, 1, $false, "67×77=", 2)
$d.Content.Find.Execute("95×96=", $true, $false, $false, $false, $false, $true, 1, $false, "88×14=", 2)
$d.Content.Find.Execute("95×36=", $true, $false, $false, $false, $false, $true, 1, $false, "25×83=", 2)
